$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 88, shifting existing rows 88:199 down to 89:200
$ws.Rows("88").Insert()

# Populate the newly inserted row 88 with the new data record
$ws.Cells.Item(88, 1).Value = 8
$ws.Cells.Item(88, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(88, 3).Value = "Coquimbo"
$ws.Cells.Item(88, 4).Value = 44638
$ws.Cells.Item(88, 5).Value = 4
$ws.Cells.Item(88, 6).Value = 100112031
$ws.Cells.Item(88, 7).Value = "Poroto verde"
$ws.Cells.Item(88, 8).Value = "Magnum"
$ws.Cells.Item(88, 9).Value = "Primera"
$ws.Cells.Item(88, 10).Value = 440
$ws.Cells.Item(88, 11).Value = 25000
$ws.Cells.Item(88, 12).Value = 26000
$ws.Cells.Item(88, 13).Value = 25500
$ws.Cells.Item(88, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(88, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(88, 16).Value = 1020
$ws.Cells.Item(88, 17).Value = 25
$ws.Cells.Item(88, 18).Value = "Hortaliza"
